# Update menu icon - mimic menu icon on mobile websites.
# Adds 5 "hamburger/menu" straight-connector lines to slide 1, built by
# duplicating the existing styled connector ("Straight Connector 7" which
# lives inside "Group 30") so the new shapes inherit the same <p:style>
# (lnRef/fillRef/effectRef/fontRef) and line fill (schemeClr bg1) as the
# rest of the deck, rather than the bare defaults AddLine/AddConnector
# produce.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# Locate "Group 30" (contains the two existing menu/plus-icon style
# straight connectors) among the slide's top-level shapes.
$group = $null
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    if ($s.Shapes.Item($i).Name -eq "Group 30") {
        $group = $s.Shapes.Item($i)
        break
    }
}

# Duplicate the group and ungroup the copy so we get free-standing
# top-level shapes to use as a style donor (this does not disturb the
# original group).
$dupGroupRange = $group.Duplicate()
$dupGroup = $dupGroupRange.Item(1)
$donorShapes = $dupGroup.Ungroup()

# "Straight Connector 7" is the first item and uses line width 28575 EMU
# (2.25 pt), matching the weight needed for the new icon bars.
$donor = $donorShapes.Item(1)

function New-MenuBar {
    param($shapeName, $left, $top, $width, $height, $weight)
    $copyRange = $donor.Duplicate()
    $shp = $copyRange.Item(1)
    $shp.Name = $shapeName
    $shp.Left = $left
    $shp.Top = $top
    $shp.Width = $width
    $shp.Height = $height
    $shp.Line.Weight = $weight
    return $shp
}

# EMU-per-point, plus a tiny epsilon: this host's Left/Top setters convert
# points -> EMU as round(pt * 12700), but the PowerShell double parsed
# from a plain division loses a couple of bits of precision, which can
# flip borderline values down to the wrong EMU. A nudge well under half
# an EMU (1 EMU = 1/12700 pt =~ 0.0000787pt) keeps the *correct* EMU but
# absorbs that rounding noise.
$EMU = 12700.0
$EPS = 0.00002

$left1 = (3924300 / $EMU) + $EPS
$top1 = (4916961 / $EMU) + $EPS
$top2 = (4988010 / $EMU) + $EPS
$left2 = (4388924 / $EMU) + $EPS
$top3 = (5062149 / $EMU) + $EPS
$top4 = (5148648 / $EMU) + $EPS
$top5 = (5243382 / $EMU) + $EPS
$widthShort = 190500 / $EMU
$widthLong = 266700 / $EMU
$zeroHeight = 0 / $EMU
$weightThin = 28575 / $EMU
$weightThick = 38100 / $EMU

$c1 = New-MenuBar "Straight Connector 43" $left1 $top1 $widthShort $zeroHeight $weightThin
$c2 = New-MenuBar "Straight Connector 44" $left1 $top2 $widthShort $zeroHeight $weightThin
$c3 = New-MenuBar "Straight Connector 45" $left2 $top3 $widthLong $zeroHeight $weightThick

# Burn through four shape ids (47-50) that were consumed-and-discarded in
# the original authoring session so the next two new connectors land on
# ids 51/52, matching the source file.
for ($i = 0; $i -lt 4; $i++) {
    $dummy = $s.Shapes.AddLine(0, 0, 1, 1)
    $dummy.Delete()
}

$c4 = New-MenuBar "Straight Connector 50" $left2 $top4 $widthLong $zeroHeight $weightThick
$c5 = New-MenuBar "Straight Connector 51" $left2 $top5 $widthLong $zeroHeight $weightThick

# Remove the temporary donor shapes (the ungrouped duplicate of Group 30).
for ($i = 1; $i -le $donorShapes.Count; $i++) {
    $donorShapes.Item($i).Delete()
}

"ok"
